$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 301.625
$ws.Range("I5").Value = 354.8
$ws.Range("J5").Value = 213
$ws.Range("K5").Value = 354.8
$ws.Range("L5").Value = 213
$ws.Range("M5").Value = -239.8
$ws.Range("N5").Value = -443
$ws.Range("H11").Value = 96
$ws.Range("I11").Value = 96
$ws.Range("K11").Value = 96
$ws.Range("M11").Value = 44
$ws.Range("H32").Value = 4667
$ws.Range("I32").Value = 2000.5
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 2000.5
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -1674.5
$ws.Range("N32").Value = -10652
$ws.Range("H40").Value = 2447.6667
$ws.Range("I40").Value = 967.4286
$ws.Range("J40").Value = 4520
$ws.Range("K40").Value = 967.4286
$ws.Range("L40").Value = 4520
$ws.Range("M40").Value = -792.4286
$ws.Range("N40").Value = -4870
$ws.Range("H53").Value = 474.53333
$ws.Range("J53").Value = 140.75
$ws.Range("L53").Value = 140.75
$ws.Range("N53").Value = -1414.75
$ws.Range("H69").Value = 3007.5
$ws.Range("J69").Value = 3007.5
$ws.Range("L69").Value = 9022.5
$ws.Range("N69").Value = -10770.5
$ws.Range("H72").Value = 3007.5
$ws.Range("J72").Value = 3007.5
$ws.Range("L72").Value = 27067.5
$ws.Range("N72").Value = -35803.5
$ws.Range("H113").Value = 11726.818
$ws.Range("I113").Value = 15599.4
$ws.Range("K113").Value = 15599.4
$ws.Range("M113").Value = -12345.4
$ws.Range("H115").Value = 449.2
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = $null
$ws.Range("H74").Value = 1867.75
$ws.Range("I74").Value = 1823.6666
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1823.6666
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -949.6666
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1867.75
$ws.Range("I77").Value = 1823.6666
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 9118.333000000001
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -4750.333000000001
$ws.Range("N77").Value = -18736
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H132").Value = 1708
$ws.Range("I132").Value = 992
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 2976
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -446
$ws.Range("N132").Value = -18060.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9399.6
$ws.Range("I86").Value = 7999.5
$ws.Range("K86").Value = 7999.5
$ws.Range("M86").Value = -6876.5
$ws.Range("H89").Value = 9399.6
$ws.Range("I89").Value = 7999.5
$ws.Range("K89").Value = 39997.5
$ws.Range("M89").Value = -34381.5
$ws.Range("H99").Value = 1644.4286
$ws.Range("I99").Value = 1478.5
$ws.Range("K99").Value = 1478.5
$ws.Range("M99").Value = 19.5
$ws.Range("H105").Value = 804.75
$ws.Range("I105").Value = 812.5714
$ws.Range("J105").Value = 750
$ws.Range("K105").Value = 812.5714
$ws.Range("L105").Value = 750
$ws.Range("M105").Value = 934.4286
$ws.Range("N105").Value = -4244
$ws.Range("H135").Value = 33333
$ws.Range("J135").Value = 29999.5
$ws.Range("L135").Value = 29999.5
$ws.Range("N135").Value = -40139.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = $null

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 200000260
$ws.Range("I4").Value = 222222500
$ws.Range("K4").Value = 666667500
$ws.Range("M4").Value = -666667388
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("H14").Value = 2001
$ws.Range("I14").Value = 2001
$ws.Range("K14").Value = 6003
$ws.Range("M14").Value = -5830
$ws.Range("H26").Value = 74.666664
$ws.Range("I26").Value = 62
$ws.Range("K26").Value = 186
$ws.Range("M26").Value = 102
$ws.Range("H55").Value = 1150
$ws.Range("I55").Value = 812.5
$ws.Range("K55").Value = 2437.5
$ws.Range("M55").Value = -2260.5
$ws.Range("H68").Value = 4145.077
$ws.Range("I68").Value = 3613
$ws.Range("J68").Value = 4996.4
$ws.Range("K68").Value = 10839
$ws.Range("L68").Value = 14989.2
$ws.Range("M68").Value = -10028
$ws.Range("N68").Value = -16611.2
$ws.Range("H71").Value = 4145.077
$ws.Range("I71").Value = 3613
$ws.Range("J71").Value = 4996.4
$ws.Range("K71").Value = 32517
$ws.Range("L71").Value = 44967.6
$ws.Range("M71").Value = -28461
$ws.Range("N71").Value = -53079.6
$ws.Range("H128").Value = 268392.34
$ws.Range("I128").Value = 268392.34
$ws.Range("K128").Value = 805177.02
$ws.Range("M128").Value = -800197.02

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2712
$ws.Range("H23").Value = 2823.5557
$ws.Range("J23").Value = 3175
$ws.Range("L23").Value = 3175
$ws.Range("N23").Value = -3621

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2054
$ws.Range("I22").Value = 1598.8
$ws.Range("J22").Value = 2433.3333
$ws.Range("K22").Value = 1598.8
$ws.Range("L22").Value = 2433.3333
$ws.Range("M22").Value = -1303.8
$ws.Range("N22").Value = -3023.3333
$ws.Range("H27").Value = 2054
$ws.Range("I27").Value = 1598.8
$ws.Range("J27").Value = 2433.3333
$ws.Range("K27").Value = 1598.8
$ws.Range("L27").Value = 2433.3333
$ws.Range("M27").Value = -1491.8
$ws.Range("N27").Value = -2647.3333
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 10000
$ws.Range("K39").Value = 10000
$ws.Range("M39").Value = -9540
$ws.Range("H46").Value = 3987.5625
$ws.Range("I46").Value = 3877
$ws.Range("J46").Value = 4466.6665
$ws.Range("K46").Value = 3877
$ws.Range("L46").Value = 4466.6665
$ws.Range("M46").Value = -3689
$ws.Range("N46").Value = -4842.6665
$ws.Range("H55").Value = 415.86667
$ws.Range("I55").Value = 143
$ws.Range("K55").Value = 143
$ws.Range("M55").Value = 30
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = $null
$ws.Range("H61").Value = 2536
$ws.Range("I61").Value = 2612.7144
$ws.Range("K61").Value = 2612.7144
$ws.Range("M61").Value = -2410.7144
$ws.Range("H92").Value = 24000
$ws.Range("J92").Value = 24000
$ws.Range("L92").Value = 24000
$ws.Range("N92").Value = -28992
$ws.Range("H100").Value = 4128.6
$ws.Range("I100").Value = 4128.6
$ws.Range("K100").Value = 4128.6
$ws.Range("M100").Value = -3587.6
$ws.Range("H101").Value = 19283.2
$ws.Range("J101").Value = 19283.2
$ws.Range("L101").Value = 19283.2
$ws.Range("N101").Value = -25773.2
$ws.Range("H113").Value = 2536
$ws.Range("I113").Value = 2612.7144
$ws.Range("K113").Value = 2612.7144
$ws.Range("M113").Value = -442.7143999999998
$ws.Range("H132").Value = 2498.5
$ws.Range("I132").Value = 2748
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 8244
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -5714
$ws.Range("N132").Value = -11058.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -36108
$ws.Range("H107").Value = 573.8182
$ws.Range("I107").Value = 263.75
$ws.Range("K107").Value = 791.25
$ws.Range("M107").Value = 1128.75
$ws.Range("H122").Value = 2639.75
$ws.Range("J122").Value = 2060.75
$ws.Range("L122").Value = 6182.25
$ws.Range("N122").Value = -11082.25
$ws.Range("H136").Value = 884.1111
$ws.Range("I136").Value = 501
$ws.Range("K136").Value = 1503
$ws.Range("M136").Value = 1047
